$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-08-14 Wednesday" "2024-08-15 Thursday"
Replace-Text "837×8=" "959×8="
Replace-Text "419×8=" "234×5="
Replace-Text "248×4=" "301×2="
Replace-Text "676×8=" "809×8="
Replace-Text "365×5=" "795×6="
Replace-Text "207×2=" "234×5="
Replace-Text "115×4=" "580×4="
Replace-Text "389×5=" "280×2="
Replace-Text "780×9=" "260×9="
Replace-Text "457×4=" "931×7="
Replace-Text "128×2=" "944×4="
Replace-Text "492×2=" "569×4="
Replace-Text "616×4=" "684×8="
Replace-Text "177×8=" "761×2="
Replace-Text "376×5=" "406×5="
Replace-Text "665×5=" "727×4="
Replace-Text "210×6=" "186×6="
Replace-Text "919×8=" "764×2="
Replace-Text "570×8=" "687×9="
Replace-Text "674×4=" "548×2="
Replace-Text "841×7=" "870×4="
Replace-Text "434×5=" "108×7="
Replace-Text "668×3=" "473×3="
Replace-Text "820×6=" "877×8="
Replace-Text "746×8=" "513×7="
